# Update the "Last Updated" timestamp on the Metadata sheet
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "29 Oct 2025, 07:01 PM"

# Add the new "distance from Dma50" sheet after "1 Month Performance"
$lastSheet = $wb.Worksheets.Item("1 Month Performance")
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "distance from Dma50"

# Header row (copy formatting/style from an existing header row)
$newSheet.Range("A1").Value = "Icon"
$newSheet.Range("B1").Value = "Stock"
$newSheet.Range("C1").Value = "Distance From Sma50"
$lastSheet.Range("A1:C1").Copy()
$newSheet.Range("A1:C1").PasteSpecial(-4122)

# Data rows
$newSheet.Cells.Item(2, 1).Value = "📈"
$newSheet.Cells.Item(2, 2).Value = "NIFTYPSUBANK"
$newSheet.Cells.Item(2, 3).Value = 10.2033
$newSheet.Cells.Item(3, 1).Value = "📈"
$newSheet.Cells.Item(3, 2).Value = "NIFTYMETAL"
$newSheet.Cells.Item(3, 3).Value = 8.6247
$newSheet.Cells.Item(4, 1).Value = "📈"
$newSheet.Cells.Item(4, 2).Value = "NIFTYOILANDGAS"
$newSheet.Cells.Item(4, 3).Value = 6.396
$newSheet.Cells.Item(5, 1).Value = "📈"
$newSheet.Cells.Item(5, 2).Value = "NIFTYCOMMODITIES"
$newSheet.Cells.Item(5, 3).Value = 5.7207
$newSheet.Cells.Item(6, 1).Value = "📈"
$newSheet.Cells.Item(6, 2).Value = "CNXINFRA"
$newSheet.Cells.Item(6, 3).Value = 5.6012
$newSheet.Cells.Item(7, 1).Value = "📈"
$newSheet.Cells.Item(7, 2).Value = "CNXREALTY"
$newSheet.Cells.Item(7, 3).Value = 5.4493
$newSheet.Cells.Item(8, 1).Value = "📈"
$newSheet.Cells.Item(8, 2).Value = "NIFTYPVTBANK"
$newSheet.Cells.Item(8, 3).Value = 5.0059
$newSheet.Cells.Item(9, 1).Value = "📈"
$newSheet.Cells.Item(9, 2).Value = "BANKNIFTY"
$newSheet.Cells.Item(9, 3).Value = 4.9192
$newSheet.Cells.Item(10, 1).Value = "📈"
$newSheet.Cells.Item(10, 2).Value = "NIFTYFINSERVICE"
$newSheet.Cells.Item(10, 3).Value = 3.9783
$newSheet.Cells.Item(11, 1).Value = "📈"
$newSheet.Cells.Item(11, 2).Value = "NIFTYMIDCAP50"
$newSheet.Cells.Item(11, 3).Value = 3.9228
$newSheet.Cells.Item(12, 1).Value = "📈"
$newSheet.Cells.Item(12, 2).Value = "NIFTY"
$newSheet.Cells.Item(12, 3).Value = 3.7191
$newSheet.Cells.Item(13, 1).Value = "📈"
$newSheet.Cells.Item(13, 2).Value = "CNXENERGY"
$newSheet.Cells.Item(13, 3).Value = 3.706
$newSheet.Cells.Item(14, 1).Value = "📈"
$newSheet.Cells.Item(14, 2).Value = "CNXMIDCAP"
$newSheet.Cells.Item(14, 3).Value = 3.6313
$newSheet.Cells.Item(15, 1).Value = "📈"
$newSheet.Cells.Item(15, 2).Value = "NIFTY200"
$newSheet.Cells.Item(15, 3).Value = 3.5915
$newSheet.Cells.Item(16, 1).Value = "📈"
$newSheet.Cells.Item(16, 2).Value = "NIFTY100"
$newSheet.Cells.Item(16, 3).Value = 3.5759
$newSheet.Cells.Item(17, 1).Value = "📈"
$newSheet.Cells.Item(17, 2).Value = "NIFTY500"
$newSheet.Cells.Item(17, 3).Value = 3.3379
$newSheet.Cells.Item(18, 1).Value = "📈"
$newSheet.Cells.Item(18, 2).Value = "CNXSMALLCAP"
$newSheet.Cells.Item(18, 3).Value = 2.9205
$newSheet.Cells.Item(19, 1).Value = "📈"
$newSheet.Cells.Item(19, 2).Value = "NIFTY50VALUE20"
$newSheet.Cells.Item(19, 3).Value = 2.893
$newSheet.Cells.Item(20, 1).Value = "📈"
$newSheet.Cells.Item(20, 2).Value = "NIFTYCPSE"
$newSheet.Cells.Item(20, 3).Value = 2.837
$newSheet.Cells.Item(21, 1).Value = "📈"
$newSheet.Cells.Item(21, 2).Value = "CNXNIFTYJUNIOR"
$newSheet.Cells.Item(21, 3).Value = 2.8291
$newSheet.Cells.Item(22, 1).Value = "📈"
$newSheet.Cells.Item(22, 2).Value = "NIFTYHEALTHCARE"
$newSheet.Cells.Item(22, 3).Value = 2.162
$newSheet.Cells.Item(23, 1).Value = "📈"
$newSheet.Cells.Item(23, 2).Value = "CNXIT"
$newSheet.Cells.Item(23, 3).Value = 2.0641
$newSheet.Cells.Item(24, 1).Value = "📈"
$newSheet.Cells.Item(24, 2).Value = "NIFTYCONSUMPTION"
$newSheet.Cells.Item(24, 3).Value = 2.0245
$newSheet.Cells.Item(25, 1).Value = "📈"
$newSheet.Cells.Item(25, 2).Value = "CNXPHARMA"
$newSheet.Cells.Item(25, 3).Value = 1.573
$newSheet.Cells.Item(26, 1).Value = "📈"
$newSheet.Cells.Item(26, 2).Value = "NIFTYAUTO"
$newSheet.Cells.Item(26, 3).Value = 1.5538
$newSheet.Cells.Item(27, 1).Value = "📈"
$newSheet.Cells.Item(27, 2).Value = "NIFTYGROWSECT15"
$newSheet.Cells.Item(27, 3).Value = 1.5325
$newSheet.Cells.Item(28, 1).Value = "📈"
$newSheet.Cells.Item(28, 2).Value = "NIFTYFMCG"
$newSheet.Cells.Item(28, 3).Value = 1.3194
$newSheet.Cells.Item(29, 1).Value = "📈"
$newSheet.Cells.Item(29, 2).Value = "NIFTYCONSURDURBL"
$newSheet.Cells.Item(29, 3).Value = 0.4031
$newSheet.Cells.Item(30, 1).Value = "📈"
$newSheet.Cells.Item(30, 2).Value = "NIFTYMEDIA"
$newSheet.Cells.Item(30, 3).Value = -1.9217

# Restore original active sheet
$meta.Activate()
